$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.618.38"
$ws.Range("E2").Value = "  -0.36%  "
$ws.Range("D3").Value = "1.595.90"
$ws.Range("E3").Value = "  -0.21%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "210.59"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("D6").Value = "0.510"
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "0.0614"
$ws.Range("E8").Value = "  -0.54%  "
$ws.Range("E9").Value = "  -0.53%  "
$ws.Range("D10").Value = "19.60"
$ws.Range("E10").Value = "  +0.33%  "
$ws.Range("E11").Value = "  +0.35%  "
$ws.Range("D12").Value = "1.819.60"
$ws.Range("E12").Value = "  -0.25%  "
$ws.Range("D13").Value = "1.602.63"
$ws.Range("E13").Value = "  +0.28%  "
$ws.Range("D14").Value = "4.04"
$ws.Range("E14").Value = "  -0.09%  "
$ws.Range("E15").Value = "  -0.30%  "
$ws.Range("D16").Value = "64.48"
$ws.Range("E16").Value = "  -1.24%  "
$ws.Range("D17").Value = "26.591.60"
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("E18").Value = "  -1.85%  "
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("D20").Value = "208.56"
$ws.Range("D21").Value = "7.06"
$ws.Range("E21").Value = "  -2.60%  "
$ws.Range("D24").Value = "8.94"
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").Value = "144.92"
$ws.Range("E25").Value = "  +1.79%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").Value = "7.11"
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("E28").Value = "  -0.62%  "
$ws.Range("E29").Value = "  -0.36%  "
$ws.Range("D30").Value = "0.0506"
$ws.Range("E30").Value = "  -3.00%  "
$ws.Range("E31").Value = "  -0.54%  "
$ws.Range("E32").Value = "  +0.08%  "
$ws.Range("E33").Value = "  -0.34%  "
$ws.Range("D34").Value = "1.282.50"
$ws.Range("E34").Value = "  -0.75%  "
$ws.Range("E35").Value = "  +0.35%  "
$ws.Range("D36").Value = "1.24"
$ws.Range("E36").Value = "  +12.70%  "
$ws.Range("E37").Value = "  -3.49%  "
$ws.Range("D38").Value = "1.48"
$ws.Range("E38").Value = "  -0.99%  "
$ws.Range("E39").Value = "  -1.84%  "
$ws.Range("D40").Value = "0.821"
$ws.Range("E40").Value = "  -0.59%  "
$ws.Range("E41").Value = "  +0.33%  "
$ws.Range("E42").Value = "  -2.00%  "
$ws.Range("E43").Value = "  -1.84%  "
$ws.Range("D44").Value = "62.84"
$ws.Range("E44").Value = "  -0.48%  "
$ws.Range("D45").Value = "1.731.74"
$ws.Range("D46").Value = "89.35"
$ws.Range("E46").Value = "  -2.07%  "
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("E48").Value = "  +2.16%  "
$ws.Range("D49").Value = "0.0514"
$ws.Range("E49").Value = "  +0.79%  "
$ws.Range("D50").Value = "7.50"
$ws.Range("E50").Value = "  +1.16%  "
$ws.Range("E51").Value = "  +0.05%  "
